$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("L_L")
$ws2 = $wb.Worksheets.Item("W_L")

# --- Add the three new rows of data (L_L sheet) ---
# Row 25: BET / PCL / FL  (stays visible after the filter below)
$ws1.Range("A25").Value = "BET"
$ws1.Range("B25").Value = "PCL"
$ws1.Range("C25").Value = "FL"
$ws1.Range("D25").Value = 1
$ws1.Range("D25").NumberFormat = "0.00"
$ws1.Range("E25").Value = 1.2129000000000001
$ws1.Range("E25").NumberFormat = "0.0000000000"
$ws1.Range("F25").Value = 18.190999999999999
$ws1.Range("F25").NumberFormat = "0.0000000000"
$ws1.Range("G25").Value = "PROP"
$ws1.Range("H25").Value = "IOTC 2005"

# Row 26: YFT / PCL / FL
$ws1.Range("A26").Value = "YFT"
$ws1.Range("B26").Value = "PCL"
$ws1.Range("C26").Value = "FL"
$ws1.Range("D26").Value = 1
$ws1.Range("D26").NumberFormat = "0.00"
$ws1.Range("E26").Value = 1.2211000000000001
$ws1.Range("E26").NumberFormat = "0.0000000000"
$ws1.Range("F26").Value = 10.733000000000001
$ws1.Range("F26").NumberFormat = "0.0000000000"
$ws1.Range("G26").Value = "PROP"
$ws1.Range("H26").Value = "IOTC 2005"

# Row 27: YFT / CKL / FL
$ws1.Range("A27").Value = "YFT"
$ws1.Range("B27").Value = "CKL"
$ws1.Range("C27").Value = "FL"
$ws1.Range("D27").Value = 1
$ws1.Range("D27").NumberFormat = "0.00"
$ws1.Range("E27").Value = 1.3693
$ws1.Range("E27").NumberFormat = "0.0000000000"
$ws1.Range("F27").Value = 21.399000000000001
$ws1.Range("F27").NumberFormat = "0.0000000000"
$ws1.Range("G27").Value = "PROP"
$ws1.Range("H27").Value = "IOTC 2005"

# Re-fit the numeric columns to their (slightly) wider new contents
$ws1.Columns.Item("E:G").AutoFit()

# --- Re-apply the AutoFilter over the extended range, filtered to BET only ---
# (drop the existing filter first so the new range A1:H27 actually "sticks")
if ($ws1.AutoFilterMode) {
    $ws1.AutoFilterMode = $false
}
$ws1.Range("A1:H27").AutoFilter(1, @("BET"))

# Keep the defined _FilterDatabase name in sync with the new filter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "L_L!_FilterDatabase") {
        $n.RefersTo = "=L_L!`$A`$1:`$H`$27"
    }
}

# --- Selection bookkeeping, matching what the author's session left behind ---
$ws2.Range("K14").Select()
$ws1.Rows.Item(26).Select()

Write-Output "edit applied"
